$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "254.56")
# are written through a Text number format so Excel stores the
# literal string instead of auto-converting it to a numeric value;
# the style is then reset to "Normal" so the cell keeps the same
# (unstyled) appearance as the rest of the data rows.

$ws.Range('D2').Value = '99.106.39'
$ws.Range('E2').Value = '  +1.36%  '
$ws.Range('D3').Value = '3.303.38'
$ws.Range('E3').Value = '  -1.19%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '254.56'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '624.20'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.65%  '
$ws.Range('E7').Value = '  +29.81%  '
$ws.Range('E8').Value = '  +5.89%  '
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('E10').Value = '  +21.33%  '
$ws.Range('D11').Value = '3.301.07'
$ws.Range('E11').Value = '  -1.24%  '
$ws.Range('E12').Value = '  +0.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '39.95'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +12.44%  '
$ws.Range('D14').Value = '98.704.79'
$ws.Range('E14').Value = '  +1.22%  '
$ws.Range('E15').Value = '  +1.55%  '
$ws.Range('D16').Value = '3.922.05'
$ws.Range('E16').Value = '  -1.02%  '
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('D18').Value = '3.303.16'
$ws.Range('E18').Value = '  -1.24%  '
$ws.Range('E19').Value = '  -4.57%  '
$ws.Range('E20').Value = '  +3.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.32'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +8.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '486.22'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('E23').Value = '  +2.30%  '
$ws.Range('E24').Value = '  -2.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.64'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '88.95'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.84%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.311'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +30.86%  '
$ws.Range('D29').Value = '3.489.08'
$ws.Range('E29').Value = '  -2.68%  '
$ws.Range('E30').Value = '  -0.10%  '
$ws.Range('E31').Value = '  +12.41%  '
$ws.Range('E32').Value = '  +2.21%  '
$ws.Range('E33').Value = '  +10.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '27.82'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.65%  '
$ws.Range('E36').Value = '  +4.93%  '
$ws.Range('E37').Value = '  -1.74%  '
$ws.Range('E38').Value = '  -2.92%  '
$ws.Range('E39').Value = '  +0.17%  '
$ws.Range('E40').Value = '  -0.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '489.12'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.70%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.63'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.94%  '
$ws.Range('E43').Value = '  -3.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.785'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.27%  '
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.12'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.66%  '
$ws.Range('E47').Value = '  +1.60%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '159.19'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.16%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.33'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +16.18%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.847'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.87%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.73'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.93%  '
